# Update the "datetime" worksheet's row-2 sample values so that all the
# America/Denver (-07:00 / MST) formatted timestamps become their UTC
# (Z / +0000 / UTC) equivalents, matching the upstream fixture refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("datetime")
$ws.Activate()

# D2: ISO8601            1989-11-09T15:17:59.123-07:00  -> 1989-11-09T15:17:59.123-00:00
$ws.Range("D2").Value = "1989-11-09T15:17:59.123-00:00"

# E2: ISO8601Z           1989-11-09T15:17:59.123Z        (unchanged content)
$ws.Range("E2").Value = "1989-11-09T15:17:59.123Z"

# F2: RFC1123            Thu, 09 Nov 1989 15:17:59 MST   -> ... UTC
$ws.Range("F2").Value = "Thu, 09 Nov 1989 15:17:59 UTC"

# G2: RFC1123Z           Thu, 09 Nov 1989 15:17:59 -0700 -> ... +0000
$ws.Range("G2").Value = "Thu, 09 Nov 1989 15:17:59 +0000"

# H2: RFC3339            1989-11-09T15:17:59-07:00       -> 1989-11-09T15:17:59Z
$ws.Range("H2").Value = "1989-11-09T15:17:59Z"

# I2: RFC3339Nano        1989-11-09T15:17:59.1234567-07:00 -> ...Z
$ws.Range("I2").Value = "1989-11-09T15:17:59.1234567Z"

# J2: RFC3339NanoZ       1989-11-09T15:17:59.1234567Z    (unchanged content)
$ws.Range("J2").Value = "1989-11-09T15:17:59.1234567Z"

# K2: RFC3339Z           1989-11-09T15:17:59Z             (unchanged content)
$ws.Range("K2").Value = "1989-11-09T15:17:59Z"

# L2: RFC8222            09 Nov 89 15:17 MST             -> ... UTC
$ws.Range("L2").Value = "09 Nov 89 15:17 UTC"

# M2: RFC8222Z - was a numeric date cell, now becomes a text value
$ws.Range("M2").Value = "09 Nov 89 15:17 +0000"

# N2: RFC850             Thursday, 09-Nov-89 15:17:59 MST -> ... UTC
$ws.Range("N2").Value = "Thursday, 09-Nov-89 15:17:59 UTC"

# O2: RubyDate           Thu Nov 09 15:17:59 -0700 1989  -> ... +0000 1989
$ws.Range("O2").Value = "Thu Nov 09 15:17:59 +0000 1989"

# T2: UnixDate           Thu Nov  9 15:17:59 MST 1989    -> ... UTC 1989
$ws.Range("T2").Value = "Thu Nov  9 15:17:59 UTC 1989"

# A2 (ANSIC) keeps the same text; rewritten only so the shared-string table
# is rebuilt cleanly alongside the other edits.
$ws.Range("A2").Value = "Thu Nov  9 15:17:59 1989"

# Move the selection to T2, matching the saved view state in the target file.
$ws.Range("T2").Select()
